$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.389927864074707
$ws.Range("B1").Value = 2.404414653778076
$ws.Range("C1").Value = 2.660238742828369
$ws.Range("D1").Value = 3.934455156326294
$ws.Range("E1").Value = 4.955169200897217
